# Apply updated cryptos list values (commit: "Updated cryptos list ... with GitHub Actions").
# All data cells in this sheet were written by openpyxl as inline/shared *strings*,
# including values that look numeric (e.g. "324.76", "1.00", "42.860.54"). Excel's COM
# Range.Value setter auto-parses plain single-dot numeric-looking text into real numbers,
# so for the Price column we prefix those values with a leading apostrophe (the normal
# Excel "force text" convention) to keep them as text, matching the source data's type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.860.54'
$ws.Range('E2').Value = '  -1.57%  '

$ws.Range('D3').Value = '2.356.37'
$ws.Range('E3').Value = '  -0.88%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').Value = '''324.76'
$ws.Range('E5').Value = '  +1.30%  '

$ws.Range('D6').Value = '''103.44'
$ws.Range('E6').Value = '  -4.20%  '

$ws.Range('D7').Value = '''0.641'
$ws.Range('E7').Value = '  +0.42%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  -1.67%  '

$ws.Range('D10').Value = '''40.25'
$ws.Range('E10').Value = '  -5.52%  '

$ws.Range('D11').Value = '''0.0926'
$ws.Range('E11').Value = '  -1.42%  '

$ws.Range('D12').Value = '''8.52'
$ws.Range('E12').Value = '  -2.19%  '

$ws.Range('D13').Value = '''1.01'
$ws.Range('E13').Value = '  -4.17%  '

$ws.Range('E14').Value = '  +0.28%  '

$ws.Range('D15').Value = '''16.22'
$ws.Range('E15').Value = '  -2.28%  '

$ws.Range('D16').Value = '2.711.68'
$ws.Range('E16').Value = '  -1.13%  '

$ws.Range('D17').Value = '2.350.17'
$ws.Range('E17').Value = '  -3.61%  '

$ws.Range('D18').Value = '42.799.74'
$ws.Range('E18').Value = '  -1.72%  '

$ws.Range('D19').Value = '''7.91'
$ws.Range('E19').Value = '  +8.99%  '

$ws.Range('E20').Value = '  -1.95%  '

$ws.Range('D21').Value = '''77.09'
$ws.Range('E21').Value = '  +2.30%  '

$ws.Range('D22').Value = '''3.68'
$ws.Range('E22').Value = '  +4.90%  '

$ws.Range('D23').Value = '''265.24'
$ws.Range('E23').Value = '  +1.56%  '

$ws.Range('E24').Value = '  -7.86%  '

$ws.Range('D25').Value = '''10.10'
$ws.Range('E25').Value = '  +9.79%  '

$ws.Range('E26').Value = '  +0.05%  '

$ws.Range('D27').Value = '''11.52'
$ws.Range('E27').Value = '  -4.15%  '

$ws.Range('D28').Value = '''23.01'
$ws.Range('E28').Value = '  +0.74%  '

$ws.Range('E29').Value = '  -1.26%  '

$ws.Range('D30').Value = '''174.97'
$ws.Range('E30').Value = '  +0.39%  '

$ws.Range('D31').Value = '''3.13'
$ws.Range('E31').Value = '  -2.65%  '

$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').Value = '''35.56'
$ws.Range('E32').Value = '  -8.68%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.0901'
$ws.Range('E33').Value = '  -2.52%  '

$ws.Range('D34').Value = '''6.22'
$ws.Range('E34').Value = '  +4.07%  '

$ws.Range('E35').Value = '  +1.49%  '

$ws.Range('D36').Value = '''0.113'
$ws.Range('E36').Value = '  +6.92%  '

$ws.Range('D37').Value = '''4.58'
$ws.Range('E37').Value = '  -7.89%  '

$ws.Range('D38').Value = '''0.0360'
$ws.Range('E38').Value = '  -3.70%  '

$ws.Range('D39').Value = '''3.81'
$ws.Range('E39').Value = '  -6.87%  '

$ws.Range('D40').Value = '''2.74'
$ws.Range('E40').Value = '  -3.33%  '

$ws.Range('D41').Value = '''0.239'
$ws.Range('E41').Value = '  +2.84%  '

$ws.Range('D42').Value = '''1.49'
$ws.Range('E42').Value = '  -1.47%  '

$ws.Range('D43').Value = '''70.55'
$ws.Range('E43').Value = '  -1.90%  '

$ws.Range('D44').Value = '''94.42'
$ws.Range('E44').Value = '  +26.02%  '

$ws.Range('D45').Value = '''120.90'
$ws.Range('E45').Value = '  +7.23%  '

$ws.Range('E46').Value = '  -0.19%  '

$ws.Range('D47').Value = '''11.95'
$ws.Range('E47').Value = '  -5.29%  '

$ws.Range('D48').Value = '''5.59'
$ws.Range('E48').Value = '  -0.61%  '

$ws.Range('D49').Value = '''9.15'
$ws.Range('E49').Value = '  -1.84%  '

$ws.Range('E50').Value = '  -3.27%  '

$ws.Range('E51').Value = '  -0.05%  '
